$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1314.4642
$ws.Range("J112").Value = 1314.4642
$ws.Range("L112").Value = 3943.3926
$ws.Range("N112").Value = -6159.392599999999

$ws.Range("H116").Value = 6500.269
$ws.Range("I116").Value = 2198.2856
$ws.Range("J116").Value = 8085.2104
$ws.Range("K116").Value = 2198.2856
$ws.Range("L116").Value = 8085.2104
$ws.Range("M116").Value = 1243.7144
$ws.Range("N116").Value = -14969.2104

$ws.Range("H129").Value = 822.43616
$ws.Range("I129").Value = 297.7143
$ws.Range("J129").Value = 864.65515
$ws.Range("K129").Value = 893.1428999999999
$ws.Range("L129").Value = 2593.96545
$ws.Range("M129").Value = 4106.8571
$ws.Range("N129").Value = -12593.96545

$ws.Range("H138").Value = 2831.6702
$ws.Range("I138").Value = 1702.9231
$ws.Range("J138").Value = 3012.8271
$ws.Range("K138").Value = 5108.7693
$ws.Range("L138").Value = 9038.481299999999
$ws.Range("M138").Value = 31.23070000000007
$ws.Range("N138").Value = -19318.4813

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1006
$ws.Range("I14").Value = 1006
$ws.Range("K14").Value = 1006
$ws.Range("M14").Value = -831

$ws.Range("H16").Value = 8399.666999999999
$ws.Range("I16").Value = 300
$ws.Range("J16").Value = 12449.5
$ws.Range("K16").Value = 300
$ws.Range("L16").Value = 12449.5
$ws.Range("M16").Value = -13
$ws.Range("N16").Value = -13023.5

$ws.Range("H32").Value = 4897.3906
$ws.Range("I32").Value = 3402.758
$ws.Range("J32").Value = 8604.08
$ws.Range("K32").Value = 3402.758
$ws.Range("L32").Value = 8604.08
$ws.Range("M32").Value = -3115.758
$ws.Range("N32").Value = -9178.08

$ws.Range("H45").Value = 2935.3333
$ws.Range("I45").Value = 3122.4
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 3122.4
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -2745.4
$ws.Range("N45").Value = -2754

$ws.Range("H74").Value = 11999.8
$ws.Range("J74").Value = 8500
$ws.Range("L74").Value = 8500
$ws.Range("N74").Value = -10248

$ws.Range("H77").Value = 11999.8
$ws.Range("J77").Value = 8500
$ws.Range("L77").Value = 42500
$ws.Range("N77").Value = -51236

$ws.Range("H122").Value = 5113.7085
$ws.Range("I122").Value = 4195.5264
$ws.Range("K122").Value = 12586.5792
$ws.Range("M122").Value = -10136.5792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 36888.777
$ws.Range("J103").Value = 36888.777
$ws.Range("L103").Value = 36888.777
$ws.Range("N103").Value = -39232.777

$ws.Range("H115").Value = 33900
$ws.Range("J115").Value = 33900
$ws.Range("L115").Value = 33900
$ws.Range("N115").Value = -37034

$ws.Range("H134").Value = 2908.3125
$ws.Range("I134").Value = 1610.1666
$ws.Range("K134").Value = 4830.4998
$ws.Range("M134").Value = -2295.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2143296.8
$ws.Range("I19").Value = 2308050.5
$ws.Range("J19").Value = 1500
$ws.Range("K19").Value = 2308050.5
$ws.Range("L19").Value = 1500
$ws.Range("M19").Value = -2307880.5
$ws.Range("N19").Value = -1840

$ws.Range("H24").Value = 2143296.8
$ws.Range("I24").Value = 2308050.5
$ws.Range("J24").Value = 1500
$ws.Range("K24").Value = 2308050.5
$ws.Range("L24").Value = 1500
$ws.Range("M24").Value = -2307880.5
$ws.Range("N24").Value = -1840

$ws.Range("H31").Value = 5961.25
$ws.Range("I31").Value = 2668.56
$ws.Range("J31").Value = 11449.066
$ws.Range("K31").Value = 2668.56
$ws.Range("L31").Value = 11449.066
$ws.Range("M31").Value = -2373.56
$ws.Range("N31").Value = -12039.066

$ws.Range("H34").Value = 5961.25
$ws.Range("I34").Value = 2668.56
$ws.Range("J34").Value = 11449.066
$ws.Range("K34").Value = 2668.56
$ws.Range("L34").Value = 11449.066
$ws.Range("M34").Value = -2466.56
$ws.Range("N34").Value = -11853.066

$ws.Range("H99").Value = 8336778
$ws.Range("I99").Value = 20001612
$ws.Range("J99").Value = 4754
$ws.Range("K99").Value = 20001612
$ws.Range("L99").Value = 4754
$ws.Range("M99").Value = -20000114
$ws.Range("N99").Value = -7750

$ws.Range("H126").Value = 8336778
$ws.Range("I126").Value = 20001612
$ws.Range("J126").Value = 4754
$ws.Range("K126").Value = 60004836
$ws.Range("L126").Value = 14262
$ws.Range("M126").Value = -60002366
$ws.Range("N126").Value = -19202

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 171.375
$ws.Range("I33").Value = 150
$ws.Range("J33").Value = 192.75
$ws.Range("K33").Value = 900
$ws.Range("L33").Value = 1156.5
$ws.Range("M33").Value = -617
$ws.Range("N33").Value = -1722.5

$ws.Range("H55").Value = 5411.6665
$ws.Range("I55").Value = 980
$ws.Range("J55").Value = 5814.5454
$ws.Range("K55").Value = 2940
$ws.Range("L55").Value = 17443.6362
$ws.Range("M55").Value = -2763
$ws.Range("N55").Value = -17797.6362

$ws.Range("H68").Value = 646
$ws.Range("I68").Value = 642
$ws.Range("J68").Value = 650
$ws.Range("K68").Value = 1926
$ws.Range("L68").Value = 1950
$ws.Range("M68").Value = -1115
$ws.Range("N68").Value = -3572

$ws.Range("H71").Value = 646
$ws.Range("I71").Value = 642
$ws.Range("J71").Value = 650
$ws.Range("K71").Value = 5778
$ws.Range("L71").Value = 5850
$ws.Range("M71").Value = -1722
$ws.Range("N71").Value = -13962

$ws.Range("H92").Value = 509.94446
$ws.Range("I92").Value = 375.2353
$ws.Range("K92").Value = 1125.7059
$ws.Range("M92").Value = 122.2941000000001

$ws.Range("H102").Value = 1791.8
$ws.Range("I102").Value = 1017.2727
$ws.Range("J102").Value = 3921.75
$ws.Range("K102").Value = 3051.8181
$ws.Range("L102").Value = 11765.25
$ws.Range("M102").Value = -617.8181
$ws.Range("N102").Value = -16633.25

$ws.Range("H113").Value = 5682414.5
$ws.Range("I113").Value = 625
$ws.Range("J113").Value = 12500562
$ws.Range("K113").Value = 1875
$ws.Range("L113").Value = 37501686
$ws.Range("M113").Value = 295
$ws.Range("N113").Value = -37506026

$ws.Range("H131").Value = 798.95
$ws.Range("I131").Value = 310
$ws.Range("J131").Value = 824.6842
$ws.Range("K131").Value = 930
$ws.Range("L131").Value = 2474.0526
$ws.Range("M131").Value = 4110
$ws.Range("N131").Value = -12554.0526

$ws.Range("H140").Value = 4251.727
$ws.Range("I140").Value = 4526.9
$ws.Range("J140").Value = 1500
$ws.Range("K140").Value = 13580.7
$ws.Range("L140").Value = 4500
$ws.Range("M140").Value = -8400.699999999999
$ws.Range("N140").Value = -14860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -10346

$ws.Range("H30").Value = 10000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = ""
$ws.Range("N30").Value = -10210

$ws.Range("H53").Value = 14996.333
$ws.Range("J53").Value = 14996.333
$ws.Range("L53").Value = 14996.333
$ws.Range("N53").Value = -16258.333

$ws.Range("H132").Value = 3861.4285
$ws.Range("I132").Value = 2506.2
$ws.Range("J132").Value = 7249.5
$ws.Range("K132").Value = 7518.599999999999
$ws.Range("L132").Value = 21748.5
$ws.Range("M132").Value = -4988.599999999999
$ws.Range("N132").Value = -26808.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6269.5
$ws.Range("I7").Value = 2098.3333
$ws.Range("J7").Value = 8057.143
$ws.Range("K7").Value = 2098.3333
$ws.Range("L7").Value = 8057.143
$ws.Range("M7").Value = -1986.3333
$ws.Range("N7").Value = -8281.143

$ws.Range("H20").Value = 4677.846
$ws.Range("J20").Value = 4986.5454
$ws.Range("L20").Value = 4986.5454
$ws.Range("N20").Value = -5438.5454

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = ""

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = ""

$ws.Range("H126").Value = 6269.5
$ws.Range("I126").Value = 2098.3333
$ws.Range("J126").Value = 8057.143
$ws.Range("K126").Value = 6294.999899999999
$ws.Range("L126").Value = 24171.429
$ws.Range("M126").Value = -3824.999899999999
$ws.Range("N126").Value = -29111.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 40374.43
$ws.Range("I23").Value = 28149.75
$ws.Range("J23").Value = 56674
$ws.Range("K23").Value = 28149.75
$ws.Range("L23").Value = 56674
$ws.Range("M23").Value = -27920.75
$ws.Range("N23").Value = -57132

$ws.Range("H25").Value = 39949.5
$ws.Range("J25").Value = 39949.5
$ws.Range("L25").Value = 39949.5
$ws.Range("N25").Value = -40535.5

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""

$ws.Range("H122").Value = 7654
$ws.Range("I122").Value = 5107
$ws.Range("J122").Value = 10519.375
$ws.Range("K122").Value = 15321
$ws.Range("L122").Value = 31558.125
$ws.Range("M122").Value = -12871
$ws.Range("N122").Value = -36458.125

$ws.Range("H136").Value = 4864.8696
$ws.Range("I136").Value = 4583.931
$ws.Range("K136").Value = 13751.793
$ws.Range("M136").Value = -11201.793
